$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting rows 45:139 down to 46:140
$ws.Rows.Item(45).Insert()

# Populate the new row 45 with the new record's data
$ws.Range("A45").Value = 10
$ws.Range("B45").Value = "Vega Modelo de Temuco"
$ws.Range("C45").Value = "La Araucanía"
$ws.Range("D45").Value = 45246
$ws.Range("D45").NumberFormat = $ws.Range("D46").NumberFormat
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100108
$ws.Range("H45").Value = "Tropicales y subtropicales"
$ws.Range("I45").Value = 100108004
$ws.Range("J45").Value = "Papaya"
$ws.Range("K45").Value = "Cultivar IV Región"
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 60
$ws.Range("N45").Value = 39000
$ws.Range("O45").Value = 39000
$ws.Range("P45").Value = 39000
$ws.Range("Q45").Value = "`$/caja 15 kilos granel"
$ws.Range("R45").Value = "Provincia del Elquí"
$ws.Range("S45").Value = 2600
$ws.Range("T45").Value = 15
